{"js": "// Replace the two-digit multiplication problems in the document, in the\n// same order they occur in the document (top-to-bottom, left-to-right).\n// Processing strictly in document order guarantees that a newly-written\n// value (e.g. the \"16\u00d735=\" written for the 19th problem) is never\n// re-matched by a later search for an earlier source value.\nconst replacements = [\n  [\"49\u00d766=\", \"20\u00d741=\"],\n  [\"36\u00d763=\", \"46\u00d780=\"],\n  [\"33\u00d745=\", \"25\u00d771=\"],\n  [\"23\u00d770=\", \"48\u00d730=\"],\n  [\"89\u00d767=\", \"58\u00d799=\"],\n  [\"84\u00d766=\", \"71\u00d732=\"],\n  [\"31\u00d718=\", \"30\u00d724=\"],\n  [\"25\u00d726=\", \"94\u00d728=\"],\n  [\"31\u00d762=\", \"94\u00d770=\"],\n  [\"50\u00d722=\", \"95\u00d761=\"],\n  [\"82\u00d796=\", \"79\u00d741=\"],\n  [\"94\u00d788=\", \"45\u00d783=\"],\n  [\"40\u00d724=\", \"83\u00d722=\"],\n  [\"97\u00d747=\", \"32\u00d781=\"],\n  [\"16\u00d735=\", \"75\u00d721=\"],\n  [\"51\u00d730=\", \"52\u00d782=\"],\n  [\"46\u00d723=\", \"94\u00d763=\"],\n  [\"67\u00d786=\", \"81\u00d789=\"],\n  [\"44\u00d780=\", \"16\u00d735=\"],\n  [\"49\u00d769=\", \"75\u00d754=\"],\n  [\"74\u00d758=\", \"28\u00d755=\"],\n  [\"86\u00d744=\", \"13\u00d788=\"],\n  [\"55\u00d757=\", \"16\u00d780=\"],\n  [\"97\u00d779=\", \"79\u00d791=\"],\n  [\"57\u00d751=\", \"61\u00d737=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  // Only the first match is the intended target; replacements are applied\n  // one at a time in document order so each search sees the still-original\n  // text for any not-yet-processed occurrence.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems throughout the document.\n# Pairs are listed in the order the source values appear in the document\n# (top-to-bottom, left-to-right). Executing the replacements in that same\n# order guarantees that a freshly written value (e.g. the \"16x35=\" that\n# becomes the 19th problem) is never re-matched by a later search for an\n# earlier, still-unprocessed source value (that search happened already).\n\n$wdReplaceOne   = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"49\u00d766=\", \"20\u00d741=\"),\n    @(\"36\u00d763=\", \"46\u00d780=\"),\n    @(\"33\u00d745=\", \"25\u00d771=\"),\n    @(\"23\u00d770=\", \"48\u00d730=\"),\n    @(\"89\u00d767=\", \"58\u00d799=\"),\n    @(\"84\u00d766=\", \"71\u00d732=\"),\n    @(\"31\u00d718=\", \"30\u00d724=\"),\n    @(\"25\u00d726=\", \"94\u00d728=\"),\n    @(\"31\u00d762=\", \"94\u00d770=\"),\n    @(\"50\u00d722=\", \"95\u00d761=\"),\n    @(\"82\u00d796=\", \"79\u00d741=\"),\n    @(\"94\u00d788=\", \"45\u00d783=\"),\n    @(\"40\u00d724=\", \"83\u00d722=\"),\n    @(\"97\u00d747=\", \"32\u00d781=\"),\n    @(\"16\u00d735=\", \"75\u00d721=\"),\n    @(\"51\u00d730=\", \"52\u00d782=\"),\n    @(\"46\u00d723=\", \"94\u00d763=\"),\n    @(\"67\u00d786=\", \"81\u00d789=\"),\n    @(\"44\u00d780=\", \"16\u00d735=\"),\n    @(\"49\u00d769=\", \"75\u00d754=\"),\n    @(\"74\u00d758=\", \"28\u00d755=\"),\n    @(\"86\u00d744=\", \"13\u00d788=\"),\n    @(\"55\u00d757=\", \"16\u00d780=\"),\n    @(\"97\u00d779=\", \"79\u00d791=\"),\n    @(\"57\u00d751=\", \"61\u00d737=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n}\n"}
